$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 200
$ws.Range("B7").Value = 250
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = -110
$ws.Range("B13").Value = -680
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 200
$ws.Range("B16").Value = -300
$ws.Range("B19").Value = -900
$ws.Range("B21").Value = -400
$ws.Range("B22").Value = -800

$ws.Range("B17").Select()
